$wb = $excel.ActiveWorkbook

# ---------- Step 1: "总计" sheet - insert new 2022-Q4 summary row ----------
$wsTotal = $wb.Worksheets.Item(1)

# Shift existing rows 2..7 down to 3..8 (bottom-up so we never overwrite unread data)
for ($r = 7; $r -ge 2; $r--) {
    $dest = $r + 1
    $wsTotal.Cells.Item($dest, 2).Value = $wsTotal.Cells.Item($r, 2).Value2
    $wsTotal.Cells.Item($dest, 3).Value = $wsTotal.Cells.Item($r, 3).Value2
    $wsTotal.Cells.Item($dest, 4).Value = $wsTotal.Cells.Item($r, 4).Value2
}

# Give the newly-used row 8 the same formatting as row 7 (index column style)
$wsTotal.Range("A7").Copy()
$wsTotal.Range("A8").PasteSpecial(-4122)

# New top row: 2022-Q4 summary
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 19
$wsTotal.Cells.Item(2, 4).Value = 6.75

# Recompute sequential index column (A) for all data rows
for ($r = 2; $r -le 8; $r++) {
    $wsTotal.Cells.Item($r, 1).Value = $r - 2
}

# ---------- Step 2: add the "2022-Q4" worksheet (fund holdings detail) ----------
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Copy($wsQ3)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Template only has 17 data rows (1 header + 16) -> extend format down to 20 rows
$newSheet.Range("A17:H17").Copy()
$newSheet.Range("A18:H20").PasteSpecial(-4122)

# Header row (already styled via the sheet copy)
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Data rows: columns B..G are textual in source data, force text storage so
# numeric-looking strings ("68.71", "012930", ...) are NOT coerced to numbers.
# row 2
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).NumberFormat = "@"
$newSheet.Cells.Item(2,2).Value = "012930"
$newSheet.Cells.Item(2,3).NumberFormat = "@"
$newSheet.Cells.Item(2,3).Value = "中庚价值先锋股票"
$newSheet.Cells.Item(2,4).NumberFormat = "@"
$newSheet.Cells.Item(2,4).Value = "68.71"
$newSheet.Cells.Item(2,5).NumberFormat = "@"
$newSheet.Cells.Item(2,5).Value = "94.78"
$newSheet.Cells.Item(2,6).NumberFormat = "@"
$newSheet.Cells.Item(2,6).Value = "4.39"
$newSheet.Cells.Item(2,7).NumberFormat = "@"
$newSheet.Cells.Item(2,7).Value = "3.0164"
$newSheet.Cells.Item(2,8).Value = 8

# row 3
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).NumberFormat = "@"
$newSheet.Cells.Item(3,2).Value = "920003"
$newSheet.Cells.Item(3,3).NumberFormat = "@"
$newSheet.Cells.Item(3,3).Value = "中金新锐股票A"
$newSheet.Cells.Item(3,4).NumberFormat = "@"
$newSheet.Cells.Item(3,4).Value = "19.30"
$newSheet.Cells.Item(3,5).NumberFormat = "@"
$newSheet.Cells.Item(3,5).Value = "91.44"
$newSheet.Cells.Item(3,6).NumberFormat = "@"
$newSheet.Cells.Item(3,6).Value = "6.08"
$newSheet.Cells.Item(3,7).NumberFormat = "@"
$newSheet.Cells.Item(3,7).Value = "1.1734"
$newSheet.Cells.Item(3,8).Value = 3

# row 4
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).NumberFormat = "@"
$newSheet.Cells.Item(4,2).Value = "011230"
$newSheet.Cells.Item(4,3).NumberFormat = "@"
$newSheet.Cells.Item(4,3).Value = "创金合信数字经济主题股票C"
$newSheet.Cells.Item(4,4).NumberFormat = "@"
$newSheet.Cells.Item(4,4).Value = "27.50"
$newSheet.Cells.Item(4,5).NumberFormat = "@"
$newSheet.Cells.Item(4,5).Value = "90.99"
$newSheet.Cells.Item(4,6).NumberFormat = "@"
$newSheet.Cells.Item(4,6).Value = "3.32"
$newSheet.Cells.Item(4,7).NumberFormat = "@"
$newSheet.Cells.Item(4,7).Value = "0.9130"
$newSheet.Cells.Item(4,8).Value = 8

# row 5
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).NumberFormat = "@"
$newSheet.Cells.Item(5,2).Value = "011229"
$newSheet.Cells.Item(5,3).NumberFormat = "@"
$newSheet.Cells.Item(5,3).Value = "创金合信数字经济主题股票A"
$newSheet.Cells.Item(5,4).NumberFormat = "@"
$newSheet.Cells.Item(5,4).Value = "20.16"
$newSheet.Cells.Item(5,5).NumberFormat = "@"
$newSheet.Cells.Item(5,5).Value = "90.99"
$newSheet.Cells.Item(5,6).NumberFormat = "@"
$newSheet.Cells.Item(5,6).Value = "3.32"
$newSheet.Cells.Item(5,7).NumberFormat = "@"
$newSheet.Cells.Item(5,7).Value = "0.6693"
$newSheet.Cells.Item(5,8).Value = 8

# row 6
$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).NumberFormat = "@"
$newSheet.Cells.Item(6,2).Value = "920923"
$newSheet.Cells.Item(6,3).NumberFormat = "@"
$newSheet.Cells.Item(6,3).Value = "中金新锐股票C"
$newSheet.Cells.Item(6,4).NumberFormat = "@"
$newSheet.Cells.Item(6,4).Value = "3.33"
$newSheet.Cells.Item(6,5).NumberFormat = "@"
$newSheet.Cells.Item(6,5).Value = "91.44"
$newSheet.Cells.Item(6,6).NumberFormat = "@"
$newSheet.Cells.Item(6,6).Value = "6.08"
$newSheet.Cells.Item(6,7).NumberFormat = "@"
$newSheet.Cells.Item(6,7).Value = "0.2025"
$newSheet.Cells.Item(6,8).Value = 3

# row 7
$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).NumberFormat = "@"
$newSheet.Cells.Item(7,2).Value = "000824"
$newSheet.Cells.Item(7,3).NumberFormat = "@"
$newSheet.Cells.Item(7,3).Value = "圆信永丰双红利灵活配置混合A"
$newSheet.Cells.Item(7,4).NumberFormat = "@"
$newSheet.Cells.Item(7,4).Value = "4.80"
$newSheet.Cells.Item(7,5).NumberFormat = "@"
$newSheet.Cells.Item(7,5).Value = "94.10"
$newSheet.Cells.Item(7,6).NumberFormat = "@"
$newSheet.Cells.Item(7,6).Value = "4.17"
$newSheet.Cells.Item(7,7).NumberFormat = "@"
$newSheet.Cells.Item(7,7).Value = "0.2002"
$newSheet.Cells.Item(7,8).Value = 5

# row 8
$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).NumberFormat = "@"
$newSheet.Cells.Item(8,2).Value = "920002"
$newSheet.Cells.Item(8,3).NumberFormat = "@"
$newSheet.Cells.Item(8,3).Value = "中金精选股票A"
$newSheet.Cells.Item(8,4).NumberFormat = "@"
$newSheet.Cells.Item(8,4).Value = "3.21"
$newSheet.Cells.Item(8,5).NumberFormat = "@"
$newSheet.Cells.Item(8,5).Value = "93.08"
$newSheet.Cells.Item(8,6).NumberFormat = "@"
$newSheet.Cells.Item(8,6).Value = "4.04"
$newSheet.Cells.Item(8,7).NumberFormat = "@"
$newSheet.Cells.Item(8,7).Value = "0.1297"
$newSheet.Cells.Item(8,8).Value = 2

# row 9
$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).NumberFormat = "@"
$newSheet.Cells.Item(9,2).Value = "550001"
$newSheet.Cells.Item(9,3).NumberFormat = "@"
$newSheet.Cells.Item(9,3).Value = "信诚四季红混合"
$newSheet.Cells.Item(9,4).NumberFormat = "@"
$newSheet.Cells.Item(9,4).Value = "4.65"
$newSheet.Cells.Item(9,5).NumberFormat = "@"
$newSheet.Cells.Item(9,5).Value = "82.79"
$newSheet.Cells.Item(9,6).NumberFormat = "@"
$newSheet.Cells.Item(9,6).Value = "2.71"
$newSheet.Cells.Item(9,7).NumberFormat = "@"
$newSheet.Cells.Item(9,7).Value = "0.1260"
$newSheet.Cells.Item(9,8).Value = 9

# row 10
$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).NumberFormat = "@"
$newSheet.Cells.Item(10,2).Value = "020015"
$newSheet.Cells.Item(10,3).NumberFormat = "@"
$newSheet.Cells.Item(10,3).Value = "国泰区位优势混合A"
$newSheet.Cells.Item(10,4).NumberFormat = "@"
$newSheet.Cells.Item(10,4).Value = "2.05"
$newSheet.Cells.Item(10,5).NumberFormat = "@"
$newSheet.Cells.Item(10,5).Value = "84.45"
$newSheet.Cells.Item(10,6).NumberFormat = "@"
$newSheet.Cells.Item(10,6).Value = "5.11"
$newSheet.Cells.Item(10,7).NumberFormat = "@"
$newSheet.Cells.Item(10,7).Value = "0.1048"
$newSheet.Cells.Item(10,8).Value = 4

# row 11
$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,2).NumberFormat = "@"
$newSheet.Cells.Item(11,2).Value = "008311"
$newSheet.Cells.Item(11,3).NumberFormat = "@"
$newSheet.Cells.Item(11,3).Value = "圆信永丰优选价值混合A"
$newSheet.Cells.Item(11,4).NumberFormat = "@"
$newSheet.Cells.Item(11,4).Value = "1.84"
$newSheet.Cells.Item(11,5).NumberFormat = "@"
$newSheet.Cells.Item(11,5).Value = "94.02"
$newSheet.Cells.Item(11,6).NumberFormat = "@"
$newSheet.Cells.Item(11,6).Value = "4.02"
$newSheet.Cells.Item(11,7).NumberFormat = "@"
$newSheet.Cells.Item(11,7).Value = "0.0740"
$newSheet.Cells.Item(11,8).Value = 7

# row 12
$newSheet.Cells.Item(12,1).Value = 10
$newSheet.Cells.Item(12,2).NumberFormat = "@"
$newSheet.Cells.Item(12,2).Value = "015594"
$newSheet.Cells.Item(12,3).NumberFormat = "@"
$newSheet.Cells.Item(12,3).Value = "国泰区位优势混合C"
$newSheet.Cells.Item(12,4).NumberFormat = "@"
$newSheet.Cells.Item(12,4).Value = "0.75"
$newSheet.Cells.Item(12,5).NumberFormat = "@"
$newSheet.Cells.Item(12,5).Value = "84.45"
$newSheet.Cells.Item(12,6).NumberFormat = "@"
$newSheet.Cells.Item(12,6).Value = "5.11"
$newSheet.Cells.Item(12,7).NumberFormat = "@"
$newSheet.Cells.Item(12,7).Value = "0.0383"
$newSheet.Cells.Item(12,8).Value = 4

# row 13
$newSheet.Cells.Item(13,1).Value = 11
$newSheet.Cells.Item(13,2).NumberFormat = "@"
$newSheet.Cells.Item(13,2).Value = "011284"
$newSheet.Cells.Item(13,3).NumberFormat = "@"
$newSheet.Cells.Item(13,3).Value = "中信保诚龙腾精选混合"
$newSheet.Cells.Item(13,4).NumberFormat = "@"
$newSheet.Cells.Item(13,4).Value = "1.08"
$newSheet.Cells.Item(13,5).NumberFormat = "@"
$newSheet.Cells.Item(13,5).Value = "83.70"
$newSheet.Cells.Item(13,6).NumberFormat = "@"
$newSheet.Cells.Item(13,6).Value = "2.71"
$newSheet.Cells.Item(13,7).NumberFormat = "@"
$newSheet.Cells.Item(13,7).Value = "0.0293"
$newSheet.Cells.Item(13,8).Value = 9

# row 14
$newSheet.Cells.Item(14,1).Value = 12
$newSheet.Cells.Item(14,2).NumberFormat = "@"
$newSheet.Cells.Item(14,2).Value = "000825"
$newSheet.Cells.Item(14,3).NumberFormat = "@"
$newSheet.Cells.Item(14,3).Value = "圆信永丰双红利灵活配置混合C"
$newSheet.Cells.Item(14,4).NumberFormat = "@"
$newSheet.Cells.Item(14,4).Value = "0.48"
$newSheet.Cells.Item(14,5).NumberFormat = "@"
$newSheet.Cells.Item(14,5).Value = "94.10"
$newSheet.Cells.Item(14,6).NumberFormat = "@"
$newSheet.Cells.Item(14,6).Value = "4.17"
$newSheet.Cells.Item(14,7).NumberFormat = "@"
$newSheet.Cells.Item(14,7).Value = "0.0200"
$newSheet.Cells.Item(14,8).Value = 5

# row 15
$newSheet.Cells.Item(15,1).Value = 13
$newSheet.Cells.Item(15,2).NumberFormat = "@"
$newSheet.Cells.Item(15,2).Value = "163818"
$newSheet.Cells.Item(15,3).NumberFormat = "@"
$newSheet.Cells.Item(15,3).Value = "中银中小盘成长混合"
$newSheet.Cells.Item(15,4).NumberFormat = "@"
$newSheet.Cells.Item(15,4).Value = "0.77"
$newSheet.Cells.Item(15,5).NumberFormat = "@"
$newSheet.Cells.Item(15,5).Value = "87.05"
$newSheet.Cells.Item(15,6).NumberFormat = "@"
$newSheet.Cells.Item(15,6).Value = "2.32"
$newSheet.Cells.Item(15,7).NumberFormat = "@"
$newSheet.Cells.Item(15,7).Value = "0.0179"
$newSheet.Cells.Item(15,8).Value = 8

# row 16
$newSheet.Cells.Item(16,1).Value = 14
$newSheet.Cells.Item(16,2).NumberFormat = "@"
$newSheet.Cells.Item(16,2).Value = "006209"
$newSheet.Cells.Item(16,3).NumberFormat = "@"
$newSheet.Cells.Item(16,3).Value = "中信保诚新蓝筹灵活配置混合"
$newSheet.Cells.Item(16,4).NumberFormat = "@"
$newSheet.Cells.Item(16,4).Value = "0.52"
$newSheet.Cells.Item(16,5).NumberFormat = "@"
$newSheet.Cells.Item(16,5).Value = "81.42"
$newSheet.Cells.Item(16,6).NumberFormat = "@"
$newSheet.Cells.Item(16,6).Value = "2.72"
$newSheet.Cells.Item(16,7).NumberFormat = "@"
$newSheet.Cells.Item(16,7).Value = "0.0141"
$newSheet.Cells.Item(16,8).Value = 8

# row 17
$newSheet.Cells.Item(17,1).Value = 15
$newSheet.Cells.Item(17,2).NumberFormat = "@"
$newSheet.Cells.Item(17,2).Value = "015201"
$newSheet.Cells.Item(17,3).NumberFormat = "@"
$newSheet.Cells.Item(17,3).Value = "创金合信动态平衡混合C"
$newSheet.Cells.Item(17,4).NumberFormat = "@"
$newSheet.Cells.Item(17,4).Value = "0.21"
$newSheet.Cells.Item(17,5).NumberFormat = "@"
$newSheet.Cells.Item(17,5).Value = "67.50"
$newSheet.Cells.Item(17,6).NumberFormat = "@"
$newSheet.Cells.Item(17,6).Value = "3.60"
$newSheet.Cells.Item(17,7).NumberFormat = "@"
$newSheet.Cells.Item(17,7).Value = "0.0076"
$newSheet.Cells.Item(17,8).Value = 5

# row 18
$newSheet.Cells.Item(18,1).Value = 16
$newSheet.Cells.Item(18,2).NumberFormat = "@"
$newSheet.Cells.Item(18,2).Value = "015200"
$newSheet.Cells.Item(18,3).NumberFormat = "@"
$newSheet.Cells.Item(18,3).Value = "创金合信动态平衡混合A"
$newSheet.Cells.Item(18,4).NumberFormat = "@"
$newSheet.Cells.Item(18,4).Value = "0.19"
$newSheet.Cells.Item(18,5).NumberFormat = "@"
$newSheet.Cells.Item(18,5).Value = "67.50"
$newSheet.Cells.Item(18,6).NumberFormat = "@"
$newSheet.Cells.Item(18,6).Value = "3.60"
$newSheet.Cells.Item(18,7).NumberFormat = "@"
$newSheet.Cells.Item(18,7).Value = "0.0068"
$newSheet.Cells.Item(18,8).Value = 5

# row 19
$newSheet.Cells.Item(19,1).Value = 17
$newSheet.Cells.Item(19,2).NumberFormat = "@"
$newSheet.Cells.Item(19,2).Value = "920922"
$newSheet.Cells.Item(19,3).NumberFormat = "@"
$newSheet.Cells.Item(19,3).Value = "中金精选股票C"
$newSheet.Cells.Item(19,4).NumberFormat = "@"
$newSheet.Cells.Item(19,4).Value = "0.12"
$newSheet.Cells.Item(19,5).NumberFormat = "@"
$newSheet.Cells.Item(19,5).Value = "93.08"
$newSheet.Cells.Item(19,6).NumberFormat = "@"
$newSheet.Cells.Item(19,6).Value = "4.04"
$newSheet.Cells.Item(19,7).NumberFormat = "@"
$newSheet.Cells.Item(19,7).Value = "0.0048"
$newSheet.Cells.Item(19,8).Value = 2

# row 20
$newSheet.Cells.Item(20,1).Value = 18
$newSheet.Cells.Item(20,2).NumberFormat = "@"
$newSheet.Cells.Item(20,2).Value = "008312"
$newSheet.Cells.Item(20,3).NumberFormat = "@"
$newSheet.Cells.Item(20,3).Value = "圆信永丰优选价值混合C"
$newSheet.Cells.Item(20,4).NumberFormat = "@"
$newSheet.Cells.Item(20,4).Value = "0.09"
$newSheet.Cells.Item(20,5).NumberFormat = "@"
$newSheet.Cells.Item(20,5).Value = "94.02"
$newSheet.Cells.Item(20,6).NumberFormat = "@"
$newSheet.Cells.Item(20,6).Value = "4.02"
$newSheet.Cells.Item(20,7).NumberFormat = "@"
$newSheet.Cells.Item(20,7).Value = "0.0036"
$newSheet.Cells.Item(20,8).Value = 7

